$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of daily-expense data (row 3): date in A3, numeric figures in B3:M3
$ws.Cells.Item(3, 1).Value = 43790
$ws.Cells.Item(3, 1).NumberFormat = "mm-dd-yy"

$rowValues = @(0, 0, 0, 0, 6.5, 0, 0, 0, 0, 1, 4, 2)
for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $ws.Cells.Item(3, $i + 2).Value = $rowValues[$i]
}

# Column A best-fits the new date column
$ws.Columns.Item(1).AutoFit()

# Selection moves to J4 after entering the new row
$ws.Range("J4").Select()
